$d = $word.ActiveDocument

# 1) Remove the _GoBack bookmark currently sitting after the "s12" paragraph text.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

function Add-EmptyParagraph {
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $r = $lastPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

function Get-LastParaRangeCollapsed {
    $p = $d.Paragraphs($d.Paragraphs.Count)
    $r = $p.Range
    $r.Collapse(0)
    return $r
}

# 2) Append new paragraphs after the last paragraph ("l3: kich co hien tai").

# Empty paragraph
Add-EmptyParagraph

# Dashes paragraph
Add-EmptyParagraph
$r = Get-LastParaRangeCollapsed
$r.Text = "--------------------------------------------------------------------------------------------------------- "

# Skeleton paragraph
Add-EmptyParagraph
$r = Get-LastParaRangeCollapsed
$r.Text = "Skeleton"

# Hyperlink paragraph 1 (atmarkcafe)
Add-EmptyParagraph
$r = Get-LastParaRangeCollapsed
$d.Hyperlinks.Add($r, "http://atmarkcafe.org/tieng-viet-viet-responsive-web-voi-skeleton/?lang=vi", "", "", "http://atmarkcafe.org/tieng-viet-viet-responsive-web-voi-skeleton/?lang=vi") | Out-Null

# Hyperlink paragraph 2 (webdesign.tutsplus.com)
Add-EmptyParagraph
$r = Get-LastParaRangeCollapsed
$d.Hyperlinks.Add($r, "https://webdesign.tutsplus.com/tutorials/building-html-page-structure-with-skeleton--cms-23253", "", "", "https://webdesign.tutsplus.com/tutorials/building-html-page-structure-with-skeleton--cms-23253") | Out-Null

# Bookmark paragraph (re-create _GoBack at the new location)
Add-EmptyParagraph
$r = Get-LastParaRangeCollapsed
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

# Trailing empty paragraph
Add-EmptyParagraph

Write-Output "Done. Paragraphs: $($d.Paragraphs.Count)"
